# Updates cryptos list values (Price / Volume(1h)) per the
# Fri May 24 11:45:14 UTC 2024 GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Some values (e.g. "85.40", "31.40") parse cleanly as a
    # plain number and Excel would silently drop the
    # insignificant trailing zero. Briefly force the cell to
    # Text format so the literal string is preserved exactly,
    # then restore the default (unstyled) cell style so no
    # stray formatting is introduced, matching the source data
    # which stores these as plain inline strings.
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.347.55"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "3.701.03"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue "D5" "596.19"
$ws.Range("E5").Value = "  -2.45%  "
Set-TextValue "D6" "165.45"
$ws.Range("E6").Value = "  -5.66%  "
$ws.Range("D7").Value = "3.696.47"
$ws.Range("E7").Value = "  -4.45%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("E11").Value = "  -4.81%  "
$ws.Range("E12").Value = "  -4.08%  "
Set-TextValue "D13" "37.71"
$ws.Range("E13").Value = "  -5.72%  "
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").Value = "4.322.04"
$ws.Range("E15").Value = "  -4.30%  "
$ws.Range("D16").Value = "3.700.84"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("D17").Value = "67.432.64"
$ws.Range("E17").Value = "  -3.69%  "
Set-TextValue "D18" "17.49"
$ws.Range("E18").Value = "  +5.32%  "
Set-TextValue "D19" "7.19"
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("E20").Value = "  -3.00%  "
Set-TextValue "D21" "487.35"
$ws.Range("E21").Value = "  -3.80%  "
Set-TextValue "D22" "9.41"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("E23").Value = "  -2.42%  "
Set-TextValue "D24" "85.40"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  -6.76%  "
Set-TextValue "D26" "0.0000139"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("E27").Value = "  -3.98%  "
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("E31").Value = "  -7.69%  "
Set-TextValue "D32" "7.63"
$ws.Range("E32").Value = "  -4.08%  "
Set-TextValue "D33" "31.40"
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("D34").Value = "3.836.09"
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("D35").Value = "3.645.60"
$ws.Range("E35").Value = "  -4.24%  "
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("E37").Value = "  -0.06%  "
Set-TextValue "D38" "0.993"
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("E39").Value = "  -6.42%  "
$ws.Range("E40").Value = "  -7.72%  "
Set-TextValue "D41" "0.322"
$ws.Range("E41").Value = "  -4.53%  "
Set-TextValue "D42" "432.29"
$ws.Range("E42").Value = "  -9.66%  "
Set-TextValue "D43" "48.55"
$ws.Range("E44").Value = "  -6.18%  "
$ws.Range("E45").Value = "  -6.36%  "
Set-TextValue "D46" "8.40"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue "D48" "40.62"
$ws.Range("E48").Value = "  -6.39%  "
Set-TextValue "D49" "142.07"
$ws.Range("E49").Value = "  +1.39%  "
$ws.Range("D50").Value = "2.750.39"
$ws.Range("E50").Value = "  -6.50%  "
$ws.Range("E51").Value = "  -4.10%  "

Write-Output "Updated 75 cells"
